$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row above row 9 (shifts existing rows 9..56 down to 10..57) ---
$ws.Rows("9:9").Insert()

# --- Populate the newly inserted row 9 with the "data_creator.py" milestone entry ---
# (values are written in this particular order so new shared-string entries land at
#  the same indices the saved workbook uses: 60..64)
$ws.Range("B9").Value = "Not a milestone - Added a python file that was used to turn data from dataset into a format that can be used for testing the database functions and data structures."
$ws.Range("C9").Value = 45776
$ws.Range("D9").Value = "main() (data_creator.py)"
$ws.Range("F9").Value = "N/A"
$ws.Range("E9").Value = "ChatGPT was used to help create the regular expression that matches if a record has the format 'UFC [event code]'`nKaggle was used to retrieve a dataset of all UFC events from 1996-2024. This is the dataset that is modified to create the testing dataset. The link for this dataset is: https://www.kaggle.com/datasets/maksbasher/ufc-complete-dataset-all-events-1996-2024?resource=download"
$ws.Range("G9").Value = "Dataset was taken from Kaggle (more information in resources used section) and converted into testing data using Python script. The output of the Python script is now in the file 'data.txt' while the original dataset is in the file 'medium_dataset.csv'."

# Match formatting used by the other milestone rows (wrap text, top-aligned; date column formatted as a date)
$ws.Range("B9:G9").WrapText = $true
$ws.Range("B9:G9").VerticalAlignment = -4160
$ws.Range("C9").NumberFormat = "m/d/yy"

# Row 9 row height (matches the source row's custom height)
$ws.Rows("9:9").RowHeight = 122.25

# A9 stays blank but carries a distinct (non-wrapping, top-aligned) style like the row it was copied above
$ws.Range("A9").VerticalAlignment = -4160
$ws.Range("A9").WrapText = $false
$ws.Range("A9").Font.Color = 0

# --- Update the selection / active cell to match the saved view state ---
$ws.Range("G6").Select()
